$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 ("#4 Complete Setting") was left half-filled (no "how"/"comment", isDone = "No").
# Finish it: mark it done and describe how it was implemented, matching the
# formatting already used by the other fully-filled rows (pink fill, border,
# centered, wrapped text).
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E11:F11").Copy() | Out-Null
$ws.Range("E12:F12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D12").Value = "Yes"
$ws.Range("E12").Value = "Get Data for ""Text Size"", ""Text Color"", and ""Background Color"" through Dialog,`nthen apply the changes to main recyclerview and layout."
$ws.Range("F12").Value = "MainActivity`nMainActivityRecyclerView`nWordsList`nSettingDialog`nactivity_main.xml`ndialog_setting.xml"

$ws.Rows.Item(12).RowHeight = 100.8

# Leave the sheet scrolled/selected where the user ended up after the edit.
$ws.Range("C13").Select() | Out-Null
